$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 830.53845
$ws.Range("I12").Value = 1243.8572
$ws.Range("J12").Value = 348.33334
$ws.Range("K12").Value = 1243.8572
$ws.Range("L12").Value = 348.33334
$ws.Range("M12").Value = -1073.8572
$ws.Range("N12").Value = -688.33334

$ws.Range("H137").Value = 5702.7144
$ws.Range("I137").Value = 6141.579
$ws.Range("J137").Value = 5340.174
$ws.Range("K137").Value = 18424.737
$ws.Range("L137").Value = 16020.522
$ws.Range("M137").Value = -15874.737
$ws.Range("N137").Value = -21120.522

$ws.Range("H138").Value = 2118.9688
$ws.Range("I138").Value = 1416.7826
$ws.Range("K138").Value = 4250.3478
$ws.Range("M138").Value = 889.6522000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 333699.5
$ws.Range("I4").Value = 333699.5
$ws.Range("K4").Value = 333699.5
$ws.Range("M4").Value = -333583.5

$ws.Range("H5").Value = 196.42857
$ws.Range("I5").Value = 195
$ws.Range("K5").Value = 195
$ws.Range("M5").Value = -83

$ws.Range("H45").Value = 11322.088
$ws.Range("I45").Value = 11278.228
$ws.Range("J45").Value = 11402.5
$ws.Range("K45").Value = 11278.228
$ws.Range("L45").Value = 11402.5
$ws.Range("M45").Value = -10901.228
$ws.Range("N45").Value = -12156.5

$ws.Range("H46").Value = 7640.6
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7640.6
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 7640.6
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -8278.6

$ws.Range("H80").Value = 18666.666

$ws.Range("H83").Value = 18666.666

$ws.Range("H132").Value = 3312.875
$ws.Range("I132").Value = 3101.647
$ws.Range("J132").Value = 4509.8335
$ws.Range("K132").Value = 9304.940999999999
$ws.Range("L132").Value = 13529.5005
$ws.Range("M132").Value = -6774.940999999999
$ws.Range("N132").Value = -18589.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 196.42857
$ws.Range("I4").Value = 195
$ws.Range("K4").Value = 195
$ws.Range("M4").Value = -80

$ws.Range("H12").Value = 3566.3333
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 5099.5
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 5099.5
$ws.Range("M12").Value = -332
$ws.Range("N12").Value = -5435.5

$ws.Range("H82").Value = 17376.166
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766

$ws.Range("H85").Value = 17376.166
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

$ws.Range("H105").Value = 9490.968000000001
$ws.Range("J105").Value = 33937
$ws.Range("L105").Value = 33937
$ws.Range("N105").Value = -37431

$ws.Range("H134").Value = 6825.32
$ws.Range("I134").Value = 5769.0557
$ws.Range("J134").Value = 9541.429
$ws.Range("K134").Value = 17307.1671
$ws.Range("L134").Value = 28624.287
$ws.Range("M134").Value = -14772.1671
$ws.Range("N134").Value = -33694.287

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 66670984
$ws.Range("I3").Value = 66670984
$ws.Range("K3").Value = 66670984
$ws.Range("M3").Value = -66670871

$ws.Range("H7").Value = 269.72726
$ws.Range("I7").Value = 181.85715
$ws.Range("J7").Value = 423.5
$ws.Range("K7").Value = 181.85715
$ws.Range("L7").Value = 423.5
$ws.Range("M7").Value = -68.85714999999999
$ws.Range("N7").Value = -649.5

$ws.Range("H31").Value = 3858.1667
$ws.Range("I31").Value = 1857.375
$ws.Range("K31").Value = 1857.375
$ws.Range("M31").Value = -1562.375

$ws.Range("H34").Value = 3858.1667
$ws.Range("I34").Value = 1857.375
$ws.Range("K34").Value = 1857.375
$ws.Range("M34").Value = -1655.375

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws.Range("H75").Value = 25260
$ws.Range("J75").Value = 25260
$ws.Range("L75").Value = 25260
$ws.Range("N75").Value = -27256

$ws.Range("H76").Value = 5706.5713
$ws.Range("I76").Value = 5706.5713
$ws.Range("K76").Value = 5706.5713
$ws.Range("M76").Value = -5391.5713

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws.Range("H78").Value = 25260
$ws.Range("J78").Value = 25260
$ws.Range("L78").Value = 75780
$ws.Range("N78").Value = -85764

$ws.Range("H79").Value = 5706.5713
$ws.Range("I79").Value = 5706.5713
$ws.Range("K79").Value = 5706.5713
$ws.Range("M79").Value = -4614.5713

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 5003.5
$ws.Range("I70").Value = 2006
$ws.Range("K70").Value = 6018
$ws.Range("M70").Value = -5703

$ws.Range("H73").Value = 5003.5
$ws.Range("I73").Value = 2006
$ws.Range("K73").Value = 6018
$ws.Range("M73").Value = -4926

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 54999
$ws.Range("J74").Value = 54999
$ws.Range("L74").Value = 54999
$ws.Range("N74").Value = -56871

$ws.Range("H77").Value = 54999
$ws.Range("J77").Value = 54999
$ws.Range("L77").Value = 164997
$ws.Range("N77").Value = -174357

$ws.Range("H80").Value = 25796.6
$ws.Range("I80").Value = 52996.5
$ws.Range("K80").Value = 52996.5
$ws.Range("M80").Value = -51998.5

$ws.Range("H83").Value = 25796.6
$ws.Range("I83").Value = 52996.5
$ws.Range("K83").Value = 264982.5
$ws.Range("M83").Value = -259990.5

$ws.Range("H113").Value = 7798.0527
$ws.Range("I113").Value = 9714.538
$ws.Range("J113").Value = 3645.6667
$ws.Range("K113").Value = 9714.538
$ws.Range("L113").Value = 3645.6667
$ws.Range("M113").Value = -7544.538
$ws.Range("N113").Value = -7985.6667

$ws.Range("H122").Value = 2610.0908
$ws.Range("I122").Value = 2467.9443
$ws.Range("K122").Value = 7403.8329
$ws.Range("M122").Value = -4953.8329

$ws.Range("H126").Value = 3128.25
$ws.Range("J126").Value = 4293.6665
$ws.Range("L126").Value = 12880.9995
$ws.Range("N126").Value = -17820.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1085.25
$ws.Range("I22").Value = 1006.6
$ws.Range("J22").Value = 1216.3334
$ws.Range("K22").Value = 1006.6
$ws.Range("L22").Value = 1216.3334
$ws.Range("M22").Value = -711.6
$ws.Range("N22").Value = -1806.3334

$ws.Range("H27").Value = 1085.25
$ws.Range("I27").Value = 1006.6
$ws.Range("J27").Value = 1216.3334
$ws.Range("K27").Value = 1006.6
$ws.Range("L27").Value = 1216.3334
$ws.Range("M27").Value = -899.6
$ws.Range("N27").Value = -1430.3334

$ws.Range("H40").Value = 7984.8184
$ws.Range("I40").Value = 8393.223
$ws.Range("J40").Value = 6147
$ws.Range("K40").Value = 8393.223
$ws.Range("L40").Value = 6147
$ws.Range("M40").Value = -8257.223
$ws.Range("N40").Value = -6419

$ws.Range("H46").Value = 1407.2
$ws.Range("I46").Value = 1490.2222
$ws.Range("J46").Value = 660
$ws.Range("K46").Value = 1490.2222
$ws.Range("L46").Value = 660
$ws.Range("M46").Value = -1302.2222
$ws.Range("N46").Value = -1036

$ws.Range("H68").Value = 44645.145
$ws.Range("I68").Value = 3129
$ws.Range("J68").Value = 100000
$ws.Range("K68").Value = 3129
$ws.Range("L68").Value = 100000
$ws.Range("M68").Value = -2380
$ws.Range("N68").Value = -101498

$ws.Range("H71").Value = 44645.145
$ws.Range("I71").Value = 3129
$ws.Range("J71").Value = 100000
$ws.Range("K71").Value = 15645
$ws.Range("L71").Value = 500000
$ws.Range("M71").Value = -11901
$ws.Range("N71").Value = -507488

$ws.Range("H76").Value = 34888.5
$ws.Range("J76").Value = 34888.5
$ws.Range("L76").Value = 34888.5
$ws.Range("N76").Value = -35564.5

$ws.Range("H79").Value = 34888.5
$ws.Range("J79").Value = 34888.5
$ws.Range("L79").Value = 34888.5
$ws.Range("N79").Value = -37228.5

$ws.Range("H99").Value = 30285
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H129").Value = 62305.668
$ws.Range("J129").Value = 62305.668
$ws.Range("L129").Value = 62305.668
$ws.Range("N129").Value = -72305.66800000001

$ws.Range("H132").Value = 4844.1377
$ws.Range("I132").Value = 4395
$ws.Range("J132").Value = 6023.125
$ws.Range("K132").Value = 13185
$ws.Range("L132").Value = 18069.375
$ws.Range("M132").Value = -10655
$ws.Range("N132").Value = -23129.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 12504750
$ws.Range("I3").Value = 16669667
$ws.Range("K3").Value = 16669667
$ws.Range("M3").Value = -16669553

$ws.Range("H11").Value = 10000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10284
